$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4766347.5
$ws.Range("I137").Value = 9097068
$ws.Range("J137").Value = 2555.3
$ws.Range("K137").Value = 27291204
$ws.Range("L137").Value = 7665.900000000001
$ws.Range("M137").Value = -27288654
$ws.Range("N137").Value = -12765.9
$ws.Range("H141").Value = 362630.25
$ws.Range("I141").Value = 1719.8636
$ws.Range("K141").Value = 5159.5908
$ws.Range("M141").Value = 20.40920000000006

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 170858.14
$ws.Range("I6").Value = 366668
$ws.Range("K6").Value = 366668
$ws.Range("M6").Value = -366495
$ws.Range("H45").Value = 1635.9445
$ws.Range("I45").Value = 1075.1724
$ws.Range("J45").Value = 3959.1428
$ws.Range("K45").Value = 1075.1724
$ws.Range("L45").Value = 3959.1428
$ws.Range("M45").Value = -698.1723999999999
$ws.Range("N45").Value = -4713.1428
$ws.Range("H61").Value = 3437.5881
$ws.Range("I61").Value = 1606.5555
$ws.Range("J61").Value = 5497.5
$ws.Range("K61").Value = 1606.5555
$ws.Range("L61").Value = 5497.5
$ws.Range("M61").Value = -1394.5555
$ws.Range("N61").Value = -5921.5
$ws.Range("H74").Value = 1295
$ws.Range("I74").Value = 1088.875
$ws.Range("J74").Value = 1707.25
$ws.Range("K74").Value = 1088.875
$ws.Range("L74").Value = 1707.25
$ws.Range("M74").Value = -214.875
$ws.Range("N74").Value = -3455.25
$ws.Range("H77").Value = 1295
$ws.Range("I77").Value = 1088.875
$ws.Range("J77").Value = 1707.25
$ws.Range("K77").Value = 5444.375
$ws.Range("L77").Value = 8536.25
$ws.Range("M77").Value = -1076.375
$ws.Range("N77").Value = -17272.25
$ws.Range("H88").Value = 2115.1428
$ws.Range("I88").Value = 1961.2
$ws.Range("J88").Value = 2500
$ws.Range("K88").Value = 1961.2
$ws.Range("L88").Value = 2500
$ws.Range("M88").Value = -1555.2
$ws.Range("N88").Value = -3312
$ws.Range("H91").Value = 2115.1428
$ws.Range("I91").Value = 1961.2
$ws.Range("J91").Value = 2500
$ws.Range("K91").Value = 1961.2
$ws.Range("L91").Value = 2500
$ws.Range("M91").Value = -557.2
$ws.Range("N91").Value = -5308
$ws.Range("H92").Value = 24079.8
$ws.Range("J92").Value = 24079.8
$ws.Range("L92").Value = 24079.8
$ws.Range("N92").Value = -29071.8
$ws.Range("H132").Value = 43484090
$ws.Range("I132").Value = 55561332
$ws.Range("J132").Value = 5998
$ws.Range("K132").Value = 166683996
$ws.Range("L132").Value = 17994
$ws.Range("M132").Value = -166681466
$ws.Range("N132").Value = -23054
$ws.Range("H133").Value = 84052.2
$ws.Range("J133").Value = 84052.2
$ws.Range("L133").Value = 84052.2
$ws.Range("N133").Value = -89112.2
$ws.Range("H136").Value = 3437.5881
$ws.Range("I136").Value = 1606.5555
$ws.Range("J136").Value = 5497.5
$ws.Range("K136").Value = 4819.666499999999
$ws.Range("L136").Value = 16492.5
$ws.Range("M136").Value = -2269.666499999999
$ws.Range("N136").Value = -21592.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 7700
$ws.Range("I11").Value = 100
$ws.Range("K11").Value = 100
$ws.Range("M11").Value = 40
$ws.Range("H80").Value = 1114.8572
$ws.Range("J80").Value = 992.2308
$ws.Range("L80").Value = 992.2308
$ws.Range("N80").Value = -2988.2308
$ws.Range("H83").Value = 1114.8572
$ws.Range("J83").Value = 992.2308
$ws.Range("L83").Value = 4961.154
$ws.Range("N83").Value = -14945.154
$ws.Range("H86").Value = 2370
$ws.Range("I86").Value = 1629.1666
$ws.Range("K86").Value = 1629.1666
$ws.Range("M86").Value = -506.1666
$ws.Range("H89").Value = 2370
$ws.Range("I89").Value = 1629.1666
$ws.Range("K89").Value = 8145.833000000001
$ws.Range("M89").Value = -2529.833000000001
$ws.Range("H134").Value = 2350.6553
$ws.Range("I134").Value = 2157.3333
$ws.Range("J134").Value = 2557.7856
$ws.Range("K134").Value = 6471.999899999999
$ws.Range("L134").Value = 7673.3568
$ws.Range("M134").Value = -3936.999899999999
$ws.Range("N134").Value = -12743.3568

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 20708.75
$ws.Range("I2").Value = 3495
$ws.Range("J2").Value = 24151.5
$ws.Range("K2").Value = 3495
$ws.Range("L2").Value = 24151.5
$ws.Range("M2").Value = -3382
$ws.Range("N2").Value = -24377.5
$ws.Range("H31").Value = 2085719.1
$ws.Range("I31").Value = 2327081.8
$ws.Range("K31").Value = 2327081.8
$ws.Range("M31").Value = -2326786.8
$ws.Range("H34").Value = 2085719.1
$ws.Range("I34").Value = 2327081.8
$ws.Range("K34").Value = 2327081.8
$ws.Range("M34").Value = -2326879.8
$ws.Range("H62").Value = 7633.3335
$ws.Range("I62").Value = 6900
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 6900
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -6276
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 7633.3335
$ws.Range("I65").Value = 6900
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 34500
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -31380
$ws.Range("N65").Value = -46240
$ws.Range("H107").Value = 3368799.2
$ws.Range("I107").Value = 5348634
$ws.Range("J107").Value = 3080.3
$ws.Range("K107").Value = 5348634
$ws.Range("L107").Value = 3080.3
$ws.Range("M107").Value = -5346714
$ws.Range("N107").Value = -6920.3
$ws.Range("H132").Value = 4449.524
$ws.Range("I132").Value = 2628
$ws.Range("J132").Value = 5360.2856
$ws.Range("K132").Value = 7884
$ws.Range("L132").Value = 16080.8568
$ws.Range("M132").Value = -5354
$ws.Range("N132").Value = -21140.8568
$ws.Range("H134").Value = 1876.129
$ws.Range("I134").Value = 1372.5927
$ws.Range("J134").Value = 5275
$ws.Range("K134").Value = 4117.7781
$ws.Range("L134").Value = 15825
$ws.Range("M134").Value = -1582.7781
$ws.Range("N134").Value = -20895

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 176
$ws.Range("I23").Value = 50
$ws.Range("J23").Value = 207.5
$ws.Range("K23").Value = 150
$ws.Range("L23").Value = 622.5
$ws.Range("M23").Value = 85
$ws.Range("N23").Value = -1092.5
$ws.Range("H34").Value = 7779.375
$ws.Range("I34").Value = 233.33333
$ws.Range("J34").Value = 9520.77
$ws.Range("K34").Value = 699.99999
$ws.Range("L34").Value = 28562.31
$ws.Range("M34").Value = -615.99999
$ws.Range("N34").Value = -28730.31
$ws.Range("H39").Value = 3487.7778
$ws.Range("I39").Value = 1166.6666
$ws.Range("J39").Value = 4648.3335
$ws.Range("K39").Value = 3499.9998
$ws.Range("L39").Value = 13945.0005
$ws.Range("M39").Value = -3205.9998
$ws.Range("N39").Value = -14533.0005
$ws.Range("H55").Value = 2440.8333
$ws.Range("I55").Value = 1060
$ws.Range("J55").Value = 3427.1428
$ws.Range("K55").Value = 3180
$ws.Range("L55").Value = 10281.4284
$ws.Range("M55").Value = -3003
$ws.Range("N55").Value = -10635.4284

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4675.0713
$ws.Range("I80").Value = 3859.182
$ws.Range("J80").Value = 7666.6665
$ws.Range("K80").Value = 3859.182
$ws.Range("L80").Value = 7666.6665
$ws.Range("M80").Value = -2861.182
$ws.Range("N80").Value = -9662.666499999999
$ws.Range("H83").Value = 4675.0713
$ws.Range("I83").Value = 3859.182
$ws.Range("J83").Value = 7666.6665
$ws.Range("K83").Value = 19295.91
$ws.Range("L83").Value = 38333.3325
$ws.Range("M83").Value = -14303.91
$ws.Range("N83").Value = -48317.3325
$ws.Range("H92").Value = 19292.166
$ws.Range("J92").Value = 19292.166
$ws.Range("L92").Value = 19292.166
$ws.Range("N92").Value = -23036.166
$ws.Range("H132").Value = 2183.7354
$ws.Range("I132").Value = 1566.7307
$ws.Range("J132").Value = 4189
$ws.Range("K132").Value = 4700.1921
$ws.Range("L132").Value = 12567
$ws.Range("M132").Value = -2170.1921
$ws.Range("N132").Value = -17627

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2030
$ws.Range("I132").Value = 1419.1818
$ws.Range("J132").Value = 3149.8333
$ws.Range("K132").Value = 4257.5454
$ws.Range("L132").Value = 9449.499899999999
$ws.Range("M132").Value = -1727.5454
$ws.Range("N132").Value = -14509.4999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 212474.56
$ws.Range("I132").Value = 304732.3
$ws.Range("J132").Value = 9507.532999999999
$ws.Range("K132").Value = 914196.8999999999
$ws.Range("L132").Value = 28522.599
$ws.Range("M132").Value = -911666.8999999999
$ws.Range("N132").Value = -33582.599
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

Write-Output "applied $(227) cell updates"